$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Time (column B) values with the new dataset
$ws.Range("B2").Value = 3.42
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2.78
$ws.Range("B5").Value = 15.18
$ws.Range("B6").Value = 24.08
$ws.Range("B7").Value = 9.78
$ws.Range("B8").Value = 2.84
$ws.Range("B9").Value = 9.04
$ws.Range("B10").Value = 10
$ws.Range("B11").Value = 17.41
$ws.Range("B12").Value = 12.98
$ws.Range("B13").Value = 15.21
$ws.Range("B14").Value = 3.68
$ws.Range("B15").Value = 30.88
$ws.Range("B16").Value = 49.83
$ws.Range("B17").Value = 6.89
$ws.Range("B18").Value = 6.94
$ws.Range("B19").Value = 3.02
$ws.Range("B20").Value = 12.72
$ws.Range("B21").Value = 23.75
$ws.Range("B22").Value = 7.48
$ws.Range("B23").Value = 15.41
$ws.Range("B24").Value = 3.63
$ws.Range("B25").Value = 16.86
$ws.Range("B26").Value = 28.61

# Move the active selection from N4 to C6
$ws.Range("C6").Select()

# Remove the duplicated (now-unused) chart defined names that were left
# over from re-creating the chart
$wb.Names.Item("_xlchart.v1.6").Delete()
$wb.Names.Item("_xlchart.v1.7").Delete()
$wb.Names.Item("_xlchart.v1.8").Delete()
